# First initial balance run for levels 1-5
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Level "base" block (rows 2-3)
$ws.Range("B2").Value = 300
$ws.Range("C2").Value = 150

# Grapevine block (rows 5-6)
$ws.Range("C5").Value = 100

# Rose Bush block (rows 8-9)
$ws.Range("C8").Value = 90

# Sunflower block (rows 11-12)
$ws.Range("C11").Value = 80

# Blue block (rows 14-15)
$ws.Range("C14").Value = 70

# Lily block (rows 17-18)
$ws.Range("C17").Value = 70

# Rose block (rows 20-21)
$ws.Range("C20").Value = 50

# Update selection to match the saved cursor position
$ws.Range("C20").Select()
